# Add data for 2022-04-13 (rename through-date from 04-04 to 04-05,
# and update the associated total counts)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (tab) to reflect the new "through" date.
$ws.Name = "Through 2022-04-05"

# Update the header label cell (shared string) that shows the same date.
$ws.Range("I1").Value = "2022 (through 04-05)"

# Update the updated monthly/total counts for the current-year column (I).
$ws.Range("I4").Value = 134
$ws.Range("I5").Value = 16
$ws.Range("I14").Value = 450
